# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals. Only the "K" column (G, rows 2-27) changes values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value
$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 3
    7  = 4
    8  = 10
    9  = 6
    10 = 4
    11 = 6
    12 = 3
    13 = 3
    14 = 2
    15 = 7
    16 = 2
    17 = 7
    18 = 0
    19 = 6
    20 = 4
    21 = 5
    22 = 3
    23 = 6
    24 = 2
    25 = 2
    26 = 2
    27 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
